$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clothing = @{
  2 = "Halter,Jumpsuit"
  3 = "Halter,Jumpsuit"
  4 = "Trunks,Caftan"
  5 = "Jumpsuit,Dress"
  6 = "Jumpsuit,Caftan"
  7 = "Caftan,Halter"
  8 = "Dress,Jumpsuit"
  9 = "Jumpsuit,Kaftan"
  10 = "Parka,Blouse"
  11 = "Blouse,Kaftan"
  12 = "Blouse,Jumpsuit"
  13 = "Jumpsuit,Blouse"
  14 = "Blouse,Halter"
  15 = "Blouse,Jumpsuit"
  16 = "Blouse,Jumpsuit"
  17 = "Blouse,Halter"
  18 = "Blouse,Halter"
  19 = "Jumpsuit,Blazer"
  20 = "Halter,Blouse"
  21 = "Dress,Trunks"
  22 = "Dress,Jumpsuit"
  23 = "Sweatpants,Trunks"
  24 = "Sweatpants,Dress"
  25 = "Gauchos,Jodhpurs"
  26 = "Trunks,Jumpsuit"
  27 = "Jumpsuit,Parka"
  28 = "Halter,Top"
  29 = "Blazer,Halter"
  30 = "Blazer,Blouse"
  31 = "Jumpsuit,Parka"
  32 = "Parka,Gauchos"
  33 = "Dress,Jodhpurs"
  34 = "Kaftan,Blouse"
  35 = "Parka,Gauchos"
  36 = "Blouse,Trunks"
  37 = "Kaftan,Jumpsuit"
  38 = "Sweatpants,Blouse"
  39 = "Trunks,Sweatpants"
  40 = "Blouse,Jumpsuit"
  41 = "Blouse,Kaftan"
  42 = "Halter,Blouse"
  43 = "Parka,Blouse"
  44 = "Jumpsuit,Blouse"
  45 = "Blazer,Top"
  46 = "Blazer,Top"
  47 = "Blouse,Jumpsuit"
  48 = "Caftan,Halter"
  49 = "Blouse,Trunks"
  50 = "Parka,Caftan"
  51 = "Jumpsuit,Kaftan"
  52 = "Blouse,Jumpsuit"
  53 = "Blouse,Kaftan"
  54 = "Blazer,Turtleneck"
  55 = "Jumpsuit,Blazer"
  56 = "Turtleneck,Jodhpurs"
  57 = "Halter,Parka"
  58 = "Halter,Dress"
  59 = "Blouse,Parka"
  60 = "Halter,Blazer"
  61 = "Blouse,Jumpsuit"
  62 = "Halter,Parka"
  63 = "Blouse,Caftan"
  64 = "Blouse,Jumpsuit"
  65 = "Blouse,Blazer"
  66 = "Jumpsuit,Blouse"
}

foreach ($row in $clothing.Keys) {
  $ws.Cells.Item($row, 7).Value = $clothing[$row]
}
